$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 108
$ws1.Range("F3").Value = 7411
$ws1.Range("F5").Value = 442
$ws1.Range("F6").Value = 3942
$ws1.Range("F7").Value = 316
$ws1.Range("F8").Value = 556
$ws1.Range("F9").Value = 272
$ws1.Range("F10").Value = 629
$ws1.Range("F11").Value = 117

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 50
$ws2.Range("F4").Value = 1

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 108
$ws4.Range("F4").Value = 7411
$ws4.Range("F5").Value = 50
$ws4.Range("F7").Value = 442
$ws4.Range("F8").Value = 3942
$ws4.Range("F9").Value = 316
$ws4.Range("F10").Value = 556
$ws4.Range("F11").Value = 272
$ws4.Range("F12").Value = 629
$ws4.Range("F13").Value = 1
$ws4.Range("F14").Value = 117
